$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.126.18'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '4.018.99'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'529.89"
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').Value = "'151.50"
$ws.Range('D7').Value = "'0.695"
$ws.Range('E7').Value = '  +11.33%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = "'0.748"
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('D11').Value = "'0.0000328"
$ws.Range('E11').Value = '  -4.31%  '
$ws.Range('D12').Value = "'47.90"
$ws.Range('E12').Value = '  +4.38%  '
$ws.Range('D13').Value = "'10.65"
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('D14').Value = '4.650.23'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '3.998.40'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').Value = "'14.09"
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = "'20.58"
$ws.Range('E17').Value = '  -3.79%  '
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('E19').Value = '  -2.58%  '
$ws.Range('D20').Value = '71.862.61'
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').Value = "'427.40"
$ws.Range('E21').Value = '  -3.32%  '
$ws.Range('D22').Value = "'98.23"
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('D23').Value = "'3.47"
$ws.Range('E23').Value = '  -3.76%  '
$ws.Range('D24').Value = "'4.19"
$ws.Range('E24').Value = '  +3.49%  '
$ws.Range('D25').Value = "'14.35"
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = "'11.22"
$ws.Range('E26').Value = '  -8.34%  '
$ws.Range('D27').Value = "'10.75"
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('D29').Value = "'36.67"
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').Value = '  +22.40%  '
$ws.Range('D31').Value = "'13.37"
$ws.Range('E31').Value = '  -1.84%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = "'7.16"
$ws.Range('E32').Value = '  +3.05%  '
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').Value = "'677.77"
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('D35').Value = "'44.62"
$ws.Range('E35').Value = '  +8.75%  '
$ws.Range('D36').Value = "'65.50"
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('D37').Value = "'0.444"
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  -4.16%  '
$ws.Range('D39').Value = '0.0₃0827'
$ws.Range('E39').Value = '  -8.66%  '
$ws.Range('D40').Value = "'3.39"
$ws.Range('E40').Value = '  -5.42%  '
$ws.Range('D41').Value = "'0.999"
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = "'0.998"
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('E44').Value = '  +2.51%  '
$ws.Range('E45').Value = '  +2.65%  '
$ws.Range('D46').Value = "'3.44"
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').Value = "'9.70"
$ws.Range('E47').Value = '  +5.10%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').Value = "'2.63"
$ws.Range('E48').Value = '  -7.48%  '
$ws.Range('D49').Value = "'2.99"
$ws.Range('E49').Value = '  -6.15%  '
$ws.Range('D50').Value = "'146.28"
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('E51').Value = '  -2.00%  '
